# Clear the "True" values in column A for the specified rows.
# These rows previously had A{row} = "True"; the diff shows them cleared
# to an empty inline string cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(35, 36, 40, 41, 42, 43, 44, 45, 46, 47, 48)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value = ""
}
